# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# The sheet currently spans A1:AC68; we extend it to A1:AF68 by appending
# three new columns: AD (Wins), AE (Losses), AF (Ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 68

# --- Header row (row 1): new header cells must carry the same header
# style ("s=1": bold, bordered, centered) as the existing header cells.
# Copy formatting from the last existing header cell (AC1) onto the new
# header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2..68): every team finished the season 81-81-0.
for ($row = 2; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 81   # column AD - Wins
    $ws.Cells.Item($row, 31).Value = 81   # column AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # column AF - Ties
}
